$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header/banner row at the top of the check-list sheet.
$ws.Rows.Item(1).Insert()

# Fill the new row with the browser-coverage note.
$ws.Range("A1").Value = "Все проверки необходимо провести для каждого из браузеров, указанных в Тест-плане"

# Make the note stand out in red (set before merging so the merged range
# ends up carrying the red-font style on every covered cell).
$ws.Range("A1").Font.Color = 255
$ws.Range("B1").Font.Color = 255
$ws.Range("C1").Font.Color = 255
$ws.Range("D1").Font.Color = 255

# Merge A1:D1 into a single banner cell.
$ws.Range("A1:D1").MergeCells = $true
